$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two more days of daily per-station charging data were appended below the
# existing table (rows 2-17 already hold 2026-02-01 .. 2026-02-08, one row
# per station, two rows per day). Clone the formatting of the last existing
# day (rows 16-17) down into the four new rows, then fill in the new values.
# Columns: A=date, B=site, C=charge amount(kwh), D=total charging income,
# E=service-fee income, F=order count.
$ws.Range("A16:F17").Copy()
$ws.Range("A18:F19").PasteSpecial(-4104)
$ws.Range("A16:F17").Copy()
$ws.Range("A20:F21").PasteSpecial(-4104)
$excel.CutCopyMode = $false

$newRows = @(
    @{ Row = 18; Date = 46062; C = 12630.67;            D = 11770.03; E = 4720.29;            F = 539 },
    @{ Row = 19; Date = 46062; C = 4800.3900000000003; D = 4336.58;  E = 1287.99;            F = 169 },
    @{ Row = 20; Date = 46063; C = 12506.17;            D = 11499.63; E = 4634.3900000000003; F = 506 },
    @{ Row = 21; Date = 46063; C = 5006.66;             D = 4481.95;  E = 1349.95;            F = 165 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}

# Scroll the view down to the newly-added rows and leave the selection where
# the author last left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("I22").Select() | Out-Null
